# Add a new "Project" worksheet at the end of the workbook (after "Agenda"),
# populate it with the Sales/Product analytics outline, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet ("Agenda") so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Project"

# --- Populate cell content -------------------------------------------------
# Note: cells are written in the same order the original author must have used
# so that the generated shared-string table lines up with the target order.
$ws.Range("B2").Value = "Schema "
$ws.Range("B3").Value = "Metadata"

$ws.Range("B5").Value = "Sales Analytics"
$ws.Range("C7").Value = "Monthly Sales Dashboard with Trends"
$ws.Range("C6").Value = "Daily Sales Summary"
$ws.Range("C8").Value = "Sales by Day of Week Analytics"
$ws.Range("C9").Value = "Hourly Sales Pattern"
$ws.Range("C10").Value = "Sales by Payment Mode"
$ws.Range("C11").Value = "Sales Vs Returns Analysis"
$ws.Range("C12").Value = "Quarterly Sales Performance"
$ws.Range("C13").Value = "Export Data for Dashboard"

$ws.Range("B15").Value = "Product Analytics"
$ws.Range("C16").Value = "Top Products Performance"
$ws.Range("C17").Value = "Category Performance Analysis"
$ws.Range("C18").Value = "Brand Performance Analysis"
$ws.Range("C19").Value = "ABC Analysis (Parento 80-20)"
$ws.Range("C20").Value = "Inventory Turnover Analysis"

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.33
$ws.Columns.Item(2).ColumnWidth = 14.33
$ws.Columns.Item(3).ColumnWidth = 30.67

# --- View state: zoom + selection + make this the active/selected tab -----
$ws.Range("F13").Select() | Out-Null
$excel.ActiveWindow.Zoom = 248
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
